$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

$wildcardChars = "!@#$%^&*()_+-=,.<>;':`"[]{}\|``~"

# --- Typography sheet: row 4 and row 6 ---
$wsTypo.Range("G4").Value = $wildcardChars
$wsTypo.Range("I4").Value = "0-9,a-z,A-Z,0x20"

$wsTypo.Range("G6").Value = $wildcardChars
$wsTypo.Range("I6").Value = "0-9,a-z,A-Z,0x20"

# --- Translation sheet: row 3 header additions ---
# (row 3 carries a row-level custom format; new cells would otherwise
# inherit that style index, but the source cells have no explicit style,
# so reset to the default "Normal" style after writing.)
$wsTrans.Range("G3").Value = "GB-DIRECTION"
$wsTrans.Range("G3").Style = "Normal"
$wsTrans.Range("H3").Value = "GB-ALIGNMENT"
$wsTrans.Range("H3").Style = "Normal"
$wsTrans.Range("I3").Value = "GB-TYPOGRAPHY"
$wsTrans.Range("I3").Style = "Normal"

# --- Translation sheet: row 9 change ---
$wsTrans.Range("C9").Value = "Small"

# --- Translation sheet: new row 12 (scroll list entry) ---
$wsTrans.Range("B12").Value = "SingleUseId10"
$wsTrans.Range("C12").Value = "Small"
$wsTrans.Range("D12").Value = "Left"
$wsTrans.Range("E12").Value = "LTR"
$wsTrans.Range("F12").Value = "INFO:"
$wsTrans.Range("G12").Value = "LTR"
$wsTrans.Range("H12").Value = "Center"
$wsTrans.Range("I12").Value = "Small"

Write-Output "edit complete"
